$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header for column B used to be "No"; it now reads "building_no".
$ws.Range("B1").Value = "building_no"

# Extend the bordered/bold header row by one (empty) cell in D1: copy the
# header formatting from C1 (bold, centered/top-aligned, boxed border),
# then trim the top/bottom edges so only the left/right borders remain -
# matching the "continuation" look used when the header strip is widened.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Borders.Item(8).LineStyle = 0
$ws.Range("D1").Borders.Item(9).LineStyle = 0
$ws.Application.CutCopyMode = $false

# Let Excel recompute the "best fit" widths for the two data columns, as
# happens automatically once the sheet is re-saved with its data intact.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("A8").Select()
